$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Worksheet data edits (Assignments/PA4/training.xlsx, Sheet1) ---
# G2:G23 -> =SUM(Cx:Fx); L2:L23 / M2:M23 get explicit formulas (same computed
# values as before); Q/R become "shared" formulas across two blocks
# (2:9 and 10:23), matching how the workbook's rows were originally built up.
for ($r = 2; $r -le 23; $r++) {
    $ws.Range("G$r").Formula = "=SUM(C$r`:F$r)"
    $ws.Range("L$r").Formula = "=(4*(B$r*B$r*4)+3*4*4)+2*(B$r*B$r*4)"
    $ws.Range("M$r").Formula = "=4*(4*B$r*B$r)+(B$r*B$r*3*4)"
}

$ws.Range("Q2:Q9").Formula   = "=0.001+2.288e-10*L2"
$ws.Range("Q10:Q23").Formula = "=0.001+2.288e-10*L10"
$ws.Range("R2:R9").Formula   = "=0.001596489+2.400611e-10*M2"
$ws.Range("R10:R23").Formula = "=0.001596489+2.400611e-10*M10"

# Restore the plain default selection on A1 (drops the stored topLeftCell
# scroll position left over from the previous session).
$ws.Range("A1").Select()

# A bit more data was added below the table -- materialize an extra
# (still-empty) row under the existing ones so the sheet grows to match.
$ws.Rows.Item(29).RowHeight = 14.25

# --- Chart tweak: data-label leader lines explicitly turned off ---
$co = $ws.ChartObjects().Item(4)
$chart = $co.Chart
$series2 = $chart.SeriesCollection().Item(2)
$dLbls = $series2.DataLabels()
$dLbls.ShowLeaderLines = $False
